$wb = $excel.ActiveWorkbook

# --- Rename the existing "Attached Functionality" sheets to "Germany" ---
$germany   = $wb.Worksheets.Item("Attached Functionality")
$germanyBB = $wb.Worksheets.Item("Attached Functionality xBB")
$germany.Name   = "Germany"
$germanyBB.Name = "GermanyxBB"

# --- Create the new "Belgium" market sheets by duplicating the Germany ones ---
# Insert right after GermanyxBB so the tab order becomes:
# Germany, GermanyxBB, Belgium, BelgiumxBB, Czech, CzechxBB
$germany.Copy($null, $germanyBB)
$belgium = $wb.Worksheets.Item(3)
$belgium.Name = "Belgium"
$belgium.Range("B2").Value = "Belgium Market"
$belgium.Range("B4").Value = "NGC-3478/T2266"

$germanyBB.Copy($null, $belgium)
$belgiumBB = $wb.Worksheets.Item(4)
$belgiumBB.Name = "BelgiumxBB"
$belgiumBB.Range("B2").Value = "Belgium Market"
$belgiumBB.Range("B4").Value = "NGC-3478/T2266"

# --- Fix up each sheet's selection / active cell ---
$czech     = $wb.Worksheets.Item("Czech")
$czechBB   = $wb.Worksheets.Item("CzechxBB")

$germany.Activate()
$germany.Range("A13").Select()

$germanyBB.Activate()
$germanyBB.Range("A13").Select()

$belgium.Activate()
$belgium.Range("B2:B4").Select()

$belgiumBB.Activate()
$belgiumBB.Range("A10").Select()

$czechBB.Activate()
$czechBB.Range("D15").Select()

# Czech tab ends up active/selected last, matching the target workbook state
$czech.Activate()
$czech.Range("C3").Select()
